$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 469
    3  = 470
    4  = 472
    5  = 474
    6  = 477
    7  = 478
    8  = 481
    9  = 483
    10 = 485
    11 = 487
    12 = 489
    13 = 492
    14 = 16
    15 = 63
    16 = 86
    17 = 104
    18 = 150
    19 = 187
    20 = 233
    21 = 257
    22 = 267
    23 = 326
    24 = 355
    25 = 393
    26 = 438
    27 = 454
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
